$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new Homework column
$ws.Range("J1").Value = "H06"

# Homework H06 scores for each student row (2-13)
$ws.Range("J2").Value = 9
$ws.Range("J3").Value = 9
$ws.Range("J4").Value = 9
$ws.Range("J5").Value = 9
$ws.Range("J6").Value = 0
$ws.Range("J7").Value = 11
$ws.Range("J8").Value = 11
$ws.Range("J9").Value = 0
$ws.Range("J10").Value = 9
$ws.Range("J11").Value = 0
$ws.Range("J12").Value = 8
$ws.Range("J13").Value = 0

# Update the active selection to match the author's final cursor position
$ws.Range("J11").Select()
